# Actualización automática 2025-08-28 15:05:09
#
# Applies the updated "agosto" (August) sales figures for advisor
# HIDALGO HIDALGO PEDRO GUSTAVO, propagating the new PIEDRA SINTERIZADA
# sales into the three report sheets:
#   - "VENTAS POR GRUPO"      (sales by product group per client)
#   - "VENTA MENSUAL"         (monthly sales per client)
#   - "CUMPLIMIENTO MENSUAL"  (budget fulfillment per product group)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: VENTAS POR GRUPO
#   Column L = PIEDRA SINTERIZADA
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# JARAMILLO CARVAJAL NICOLAS ESTEBAN
$wsGrupo.Range("L11").Value = 147.81

# TRUJILLO TORRES VINICIO RUBEN
$wsGrupo.Range("L21").Value = 873.8

# Summary row: count of clients with sales in each group
$wsGrupo.Range("L23").Value = "2 de 21"

# ---------------------------------------------------------------------
# Sheet: VENTA MENSUAL
#   Column F = agosto
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# JARAMILLO CARVAJAL NICOLAS ESTEBAN
$wsMensual.Range("F11").Value = 3463.44

# TRUJILLO TORRES VINICIO RUBEN
$wsMensual.Range("F21").Value = 873.8

# TOTAL row
$wsMensual.Range("F23").Value = 13969.18

# ---------------------------------------------------------------------
# Sheet: CUMPLIMIENTO MENSUAL
#   Row 15 = PIEDRA SINTERIZADA group, Row 19 = TOTAL
#   C = PRESUPUESTO, D = VENTA, E = POR CUMPLIR, F = CUMPLIMIENTO
# ---------------------------------------------------------------------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# PIEDRA SINTERIZADA row
$wsCumplimiento.Range("D15").Value = 1021.61
$wsCumplimiento.Range("E15").Value = 6443.39
$wsCumplimiento.Range("F15").Value = 0.1368533154722036

# TOTAL row
$wsCumplimiento.Range("D19").Value = 15684.49
$wsCumplimiento.Range("E19").Value = 43703.73762291769
$wsCumplimiento.Range("F19").Value = 0.2641009948905668
